$wb = $excel.ActiveWorkbook

# Sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1213.7333
$ws.Range("I40").Value = 1225.5
$ws.Range("J40").Value = 1166.6666
$ws.Range("K40").Value = 1225.5
$ws.Range("L40").Value = 1166.6666
$ws.Range("M40").Value = -1050.5
$ws.Range("N40").Value = -1516.6666

# Sheet ALC, row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 41469.4
$ws.Range("I86").Value = 1540.6
$ws.Range("J86").Value = 81398.2
$ws.Range("K86").Value = 1540.6
$ws.Range("L86").Value = 81398.2
$ws.Range("M86").Value = -417.5999999999999
$ws.Range("N86").Value = -83644.2

# Sheet ALC, row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 41469.4
$ws.Range("I89").Value = 1540.6
$ws.Range("J89").Value = 81398.2
$ws.Range("K89").Value = 7703
$ws.Range("L89").Value = 406991
$ws.Range("M89").Value = -2087
$ws.Range("N89").Value = -418223

# Sheet ALC, row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1144.575
$ws.Range("J112").Value = 1150.9429
$ws.Range("L112").Value = 3452.8287
$ws.Range("N112").Value = -5668.8287

# Sheet ALC, row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 8142.5713
$ws.Range("I113").Value = 2177.2222
$ws.Range("J113").Value = 12616.583
$ws.Range("K113").Value = 2177.2222
$ws.Range("L113").Value = 12616.583
$ws.Range("M113").Value = 1076.7778
$ws.Range("N113").Value = -19124.583

# Sheet ALC, row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 4481.62
$ws.Range("I131").Value = 713.5
$ws.Range("J131").Value = 4995.4546
$ws.Range("K131").Value = 2140.5
$ws.Range("L131").Value = 14986.3638
$ws.Range("M131").Value = 2899.5
$ws.Range("N131").Value = -25066.3638

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3828.375
$ws.Range("I32").Value = 2358.2273
$ws.Range("K32").Value = 2358.2273
$ws.Range("M32").Value = -2071.2273

# Sheet ARM, row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 8454.6
$ws.Range("I88").Value = 11531.2
$ws.Range("J88").Value = 2301.4
$ws.Range("K88").Value = 11531.2
$ws.Range("L88").Value = 2301.4
$ws.Range("M88").Value = -11125.2
$ws.Range("N88").Value = -3113.4

# Sheet ARM, row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 8454.6
$ws.Range("I91").Value = 11531.2
$ws.Range("J91").Value = 2301.4
$ws.Range("K91").Value = 11531.2
$ws.Range("L91").Value = 2301.4
$ws.Range("M91").Value = -10127.2
$ws.Range("N91").Value = -5109.4

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 8549380
$ws.Range("I122").Value = 13335454
$ws.Range("J122").Value = 2820.8572
$ws.Range("K122").Value = 40006362
$ws.Range("L122").Value = 8462.571599999999
$ws.Range("M122").Value = -40003912
$ws.Range("N122").Value = -13362.5716

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2170.9048
$ws.Range("I132").Value = 997.9666999999999
$ws.Range("J132").Value = 5103.25
$ws.Range("K132").Value = 2993.9001
$ws.Range("L132").Value = 15309.75
$ws.Range("M132").Value = -463.9000999999998
$ws.Range("N132").Value = -20369.75

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3168.1758
$ws.Range("I31").Value = 682.9474
$ws.Range("J31").Value = 7334.5884
$ws.Range("K31").Value = 682.9474
$ws.Range("L31").Value = 7334.5884
$ws.Range("M31").Value = -387.9474
$ws.Range("N31").Value = -7924.5884

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3168.1758
$ws.Range("I34").Value = 682.9474
$ws.Range("J34").Value = 7334.5884
$ws.Range("K34").Value = 682.9474
$ws.Range("L34").Value = 7334.5884
$ws.Range("M34").Value = -480.9474
$ws.Range("N34").Value = -7738.5884

# Sheet CRP, row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 72908.5
$ws.Range("I105").Value = 72908.5
$ws.Range("K105").Value = 72908.5
$ws.Range("M105").Value = -71161.5

# Sheet CUL, row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 2255.7778
$ws.Range("I46").Value = 151
$ws.Range("J46").Value = 2857.1428
$ws.Range("K46").Value = 453
$ws.Range("L46").Value = 8571.428400000001
$ws.Range("M46").Value = -362
$ws.Range("N46").Value = -8753.428400000001

# Sheet CUL, row 48
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 3600
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 3600
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 10800
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -11300

# Sheet CUL, row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3660
$ws.Range("I55").Value = 2980
$ws.Range("J55").Value = 4000
$ws.Range("K55").Value = 8940
$ws.Range("L55").Value = 12000
$ws.Range("M55").Value = -8763
$ws.Range("N55").Value = -12354

# Sheet CUL, row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 535093.4
$ws.Range("I68").Value = 1344082.6
$ws.Range("J68").Value = 1504.7021
$ws.Range("K68").Value = 4032247.8
$ws.Range("L68").Value = 4514.106299999999
$ws.Range("M68").Value = -4031436.8
$ws.Range("N68").Value = -6136.106299999999

# Sheet CUL, row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 535093.4
$ws.Range("I71").Value = 1344082.6
$ws.Range("J71").Value = 1504.7021
$ws.Range("K71").Value = 12096743.4
$ws.Range("L71").Value = 13542.3189
$ws.Range("M71").Value = -12092687.4
$ws.Range("N71").Value = -21654.3189

# Sheet CUL, row 99
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

# Sheet GSM, row 3
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1670317.1
$ws.Range("I3").Value = 3337001
$ws.Range("J3").Value = 3633.3333
$ws.Range("K3").Value = 3337001
$ws.Range("L3").Value = 3633.3333
$ws.Range("M3").Value = -3336885
$ws.Range("N3").Value = -3865.3333

# Sheet GSM, row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1274.2
$ws.Range("I113").Value = 1289.8572
$ws.Range("K113").Value = 1289.8572
$ws.Range("M113").Value = 880.1428000000001

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3053.1924
$ws.Range("I132").Value = 1835.6875
$ws.Range("J132").Value = 5001.2
$ws.Range("K132").Value = 5507.0625
$ws.Range("L132").Value = 15003.6
$ws.Range("M132").Value = -2977.0625
$ws.Range("N132").Value = -20063.6

# Sheet LTW, row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 146714.28
$ws.Range("I7").Value = 253750
$ws.Range("K7").Value = 253750
$ws.Range("M7").Value = -253638

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 33493.883
$ws.Range("I40").Value = 52780
$ws.Range("J40").Value = 5942.2856
$ws.Range("K40").Value = 52780
$ws.Range("L40").Value = 5942.2856
$ws.Range("M40").Value = -52644
$ws.Range("N40").Value = -6214.2856

# Sheet LTW, row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 146714.28
$ws.Range("I126").Value = 253750
$ws.Range("K126").Value = 761250
$ws.Range("M126").Value = -758780

# Sheet WVR, row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 69169.664
$ws.Range("J2").Value = 69169.664
$ws.Range("L2").Value = 69169.664
$ws.Range("N2").Value = -69393.664

# Sheet WVR, row 4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 49980
$ws.Range("I4").Value = 50000
$ws.Range("J4").Value = 49975
$ws.Range("K4").Value = 50000
$ws.Range("L4").Value = 49975
$ws.Range("M4").Value = -49887
$ws.Range("N4").Value = -50201

# Sheet WVR, row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3000.6667
$ws.Range("I62").Value = 3000.6667
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3000.6667
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2376.6667
$ws.Range("N62").ClearContents()

# Sheet WVR, row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 3000.6667
$ws.Range("I65").Value = 3000.6667
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 15003.3335
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -11883.3335
$ws.Range("N65").ClearContents()

# Sheet WVR, row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1491
$ws.Range("I96").Value = 1235.1666
$ws.Range("J96").Value = 2002.6666
$ws.Range("K96").Value = 1235.1666
$ws.Range("L96").Value = 2002.6666
$ws.Range("M96").Value = 137.8334
$ws.Range("N96").Value = -4748.6666

# Sheet WVR, row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 67318
$ws.Range("I100").Value = 605.8182
$ws.Range("J100").Value = 250776.5
$ws.Range("K100").Value = 1211.6364
$ws.Range("L100").Value = 501553
$ws.Range("M100").Value = -670.6364000000001
$ws.Range("N100").Value = -502635

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 47535.453
$ws.Range("I122").Value = 68525.664
$ws.Range("J122").Value = 2556.4285
$ws.Range("K122").Value = 205576.992
$ws.Range("L122").Value = 7669.2855
$ws.Range("M122").Value = -203126.992
$ws.Range("N122").Value = -12569.2855
